$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (row 6) gains a new first column "Cliente" in C6, pushing the
# previous "Número de Documento" text out of use. Replace C6's text.
$ws.Range("C6").Value = "Cliente"

# Column C needs to widen to fit the new "Cliente" header/content.
$ws.Columns("C").ColumnWidth = 23.45

# Row 6 no longer needs the explicit 30pt height - let it size back to default.
$ws.Rows(6).AutoFit()

# Move/collapse the active selection to a single cell, D4.
$ws.Range("D4").Select()
